# BOM.xlsx update
#  - Fix typo in the I2C-Wandler hint (G9): "ausgeählt" -> "ausgwählt"
#  - Specify the RGB-LED package size in D12: "Common Anode RGB LED" -> "Common Anode RGB LED (5mm)"
#  - Add a note about cable length/colour choice for the Stromkabel row (G20)
#  - Add a new "Hinweise:" section below the table (rows 25-26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Typo fix on the LCD / I2C-Wandler remark
$ws.Range("G9").Value = "Es kann auch direkt ausgwählt werden, dass ein I2C-Wandler vorinstalliert werden soll. Dies spart Arbeit und Versandkosten."

# 2) Specify LED package size
$ws.Range("D12").Value = "Common Anode RGB LED (5mm)"

# 3) New note for the USB power cable row
$ws.Range("G20").Value = "Es gibt verschiende Längen und unterschiedliche Kabelfarben zur Auswahl, solange der Anschluss passt, ist diese Auswahl egal."

# 4) New "Hinweise" block under the table
$ws.Range("B25").Value = "Hinweise:"
$ws.Range("C25").Value = "Allgemein lohnt es sich, die Teile zu Vergleichen und zu schauen, wer zur Zeit die besten Teilpreise liefert."
$ws.Range("C26").Value = "Auch Sammelbestellungen lohnen sich sehr, so schafft man niedrigere Gesamtpreise und spart Versandkosten und CO2-Emissionen."

# 5) Move the active selection like the author's saved state
$ws.Range("I44").Select()
